$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header fields (row 21) ---
$ws.Range("C21").Value = 7
$ws.Range("F21").Value = 44746
$ws.Range("I21").Value = 44772

# --- Unmerge the rows we are about to restructure ---
$ws.Range("C25:G25").UnMerge()
$ws.Range("C26:G26").UnMerge()
$ws.Range("H27:I27").UnMerge()
$ws.Range("H28:I28").UnMerge()
$ws.Range("H29:I29").UnMerge()

# --- Row 24 becomes the single "Pan" line item ---
# Force these as text (matching the original inline-string cells) so values
# like "15.00", "9%" and "45.00" aren't silently reinterpreted as numbers.
$ws.Range("H24:K24").NumberFormat = "@"
$ws.Range("C24").Value = "Pan"
$ws.Range("H24").Value = "3"
$ws.Range("I24").Value = "15.00"
$ws.Range("J24").Value = "9%"
$ws.Range("K24").Value = "45.00"

# --- Clear old product rows 25 and 26 ---
$ws.Range("C25:K26").ClearContents()

# --- Delete the now-unused rows 28 and 29 (old BTW/Totaal rows) ---
$ws.Range("A28:A29").EntireRow.Delete()

# --- Row 25: Subtotaal ---
$ws.Range("H25").Value = "Subtotaal"
$ws.Range("J25").Formula = "=J27-J26"

# --- Row 26: BTW ---
$ws.Range("H26").Value = "BTW"
$ws.Range("J26").Formula = "=J27*0.09"

# --- Row 27: Totaal ---
$ws.Range("H27").Value = "Totaal"
$ws.Range("J27").Value = 45

# --- Re-merge cells in their new layout ---
$ws.Range("H25:I25").Merge()
$ws.Range("H26:I26").Merge()
$ws.Range("H27:I27").Merge()
